{"js": "// Replace the 100 two-digit multiplication problems (20 rows x 5 cols)\n// in the worksheet table with the new set of problems from the commit.\n// Table cell text is replaced in place (table.values = ...) which keeps\n// each cell's existing paragraph/run formatting (rFonts, sz, jc) intact \u2014\n// only the <w:t> text content changes, matching the OOXML diff exactly.\n\nconst newValues = [\n  [\"39\u00d749=\", \"48\u00d727=\", \"10\u00d748=\", \"58\u00d726=\", \"30\u00d759=\"],\n  [\"29\u00d7100=\", \"40\u00d760=\", \"69\u00d743=\", \"47\u00d746=\", \"49\u00d750=\"],\n  [\"34\u00d715=\", \"43\u00d769=\", \"12\u00d761=\", \"78\u00d736=\", \"59\u00d752=\"],\n  [\"41\u00d765=\", \"76\u00d746=\", \"51\u00d711=\", \"13\u00d797=\", \"56\u00d757=\"],\n  [\"88\u00d748=\", \"29\u00d740=\", \"30\u00d711=\", \"18\u00d725=\", \"61\u00d765=\"],\n  [\"81\u00d765=\", \"75\u00d790=\", \"49\u00d730=\", \"54\u00d745=\", \"17\u00d792=\"],\n  [\"23\u00d757=\", \"47\u00d782=\", \"100\u00d743=\", \"22\u00d771=\", \"70\u00d777=\"],\n  [\"34\u00d742=\", \"36\u00d760=\", \"93\u00d768=\", \"95\u00d763=\", \"16\u00d795=\"],\n  [\"76\u00d765=\", \"68\u00d753=\", \"61\u00d750=\", \"39\u00d771=\", \"11\u00d757=\"],\n  [\"10\u00d750=\", \"69\u00d798=\", \"48\u00d772=\", \"72\u00d717=\", \"85\u00d711=\"],\n  [\"17\u00d748=\", \"38\u00d739=\", \"82\u00d759=\", \"29\u00d776=\", \"30\u00d781=\"],\n  [\"57\u00d750=\", \"11\u00d783=\", \"48\u00d787=\", \"82\u00d734=\", \"95\u00d790=\"],\n  [\"77\u00d728=\", \"50\u00d728=\", \"98\u00d761=\", \"12\u00d798=\", \"64\u00d758=\"],\n  [\"84\u00d790=\", \"20\u00d778=\", \"13\u00d789=\", \"47\u00d728=\", \"86\u00d792=\"],\n  [\"63\u00d790=\", \"44\u00d736=\", \"22\u00d719=\", \"41\u00d751=\", \"41\u00d771=\"],\n  [\"19\u00d795=\", \"54\u00d711=\", \"87\u00d776=\", \"90\u00d768=\", \"80\u00d749=\"],\n  [\"35\u00d759=\", \"19\u00d734=\", \"66\u00d783=\", \"55\u00d772=\", \"11\u00d768=\"],\n  [\"99\u00d778=\", \"46\u00d755=\", \"57\u00d780=\", \"98\u00d755=\", \"85\u00d724=\"],\n  [\"30\u00d732=\", \"47\u00d768=\", \"63\u00d721=\", \"91\u00d734=\", \"65\u00d791=\"],\n  [\"95\u00d780=\", \"12\u00d775=\", \"89\u00d797=\", \"52\u00d714=\", \"36\u00d731=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table of multiplication problems, found none.\");\n}\n\nconst table = tables.items[0];\n\n// Assigning the 2-D values grid in one shot updates each cell's existing\n// paragraph/run in place (formatting such as rFonts/sz/jc is preserved) \u2014\n// only the <w:t> text content changes, matching the OOXML diff exactly.\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# Replace the 100 two-digit multiplication problems (20 rows x 5 cols)\n# in the worksheet table with the new set of problems from the commit.\n# Setting Cell.Range.Text keeps the cell's existing paragraph/run\n# formatting (rFonts, sz, jc) intact -- only the <w:t> text changes,\n# matching the OOXML diff exactly.\n\n$newValues = @(\n  @(\"39\u00d749=\", \"48\u00d727=\", \"10\u00d748=\", \"58\u00d726=\", \"30\u00d759=\"),\n  @(\"29\u00d7100=\", \"40\u00d760=\", \"69\u00d743=\", \"47\u00d746=\", \"49\u00d750=\"),\n  @(\"34\u00d715=\", \"43\u00d769=\", \"12\u00d761=\", \"78\u00d736=\", \"59\u00d752=\"),\n  @(\"41\u00d765=\", \"76\u00d746=\", \"51\u00d711=\", \"13\u00d797=\", \"56\u00d757=\"),\n  @(\"88\u00d748=\", \"29\u00d740=\", \"30\u00d711=\", \"18\u00d725=\", \"61\u00d765=\"),\n  @(\"81\u00d765=\", \"75\u00d790=\", \"49\u00d730=\", \"54\u00d745=\", \"17\u00d792=\"),\n  @(\"23\u00d757=\", \"47\u00d782=\", \"100\u00d743=\", \"22\u00d771=\", \"70\u00d777=\"),\n  @(\"34\u00d742=\", \"36\u00d760=\", \"93\u00d768=\", \"95\u00d763=\", \"16\u00d795=\"),\n  @(\"76\u00d765=\", \"68\u00d753=\", \"61\u00d750=\", \"39\u00d771=\", \"11\u00d757=\"),\n  @(\"10\u00d750=\", \"69\u00d798=\", \"48\u00d772=\", \"72\u00d717=\", \"85\u00d711=\"),\n  @(\"17\u00d748=\", \"38\u00d739=\", \"82\u00d759=\", \"29\u00d776=\", \"30\u00d781=\"),\n  @(\"57\u00d750=\", \"11\u00d783=\", \"48\u00d787=\", \"82\u00d734=\", \"95\u00d790=\"),\n  @(\"77\u00d728=\", \"50\u00d728=\", \"98\u00d761=\", \"12\u00d798=\", \"64\u00d758=\"),\n  @(\"84\u00d790=\", \"20\u00d778=\", \"13\u00d789=\", \"47\u00d728=\", \"86\u00d792=\"),\n  @(\"63\u00d790=\", \"44\u00d736=\", \"22\u00d719=\", \"41\u00d751=\", \"41\u00d771=\"),\n  @(\"19\u00d795=\", \"54\u00d711=\", \"87\u00d776=\", \"90\u00d768=\", \"80\u00d749=\"),\n  @(\"35\u00d759=\", \"19\u00d734=\", \"66\u00d783=\", \"55\u00d772=\", \"11\u00d768=\"),\n  @(\"99\u00d778=\", \"46\u00d755=\", \"57\u00d780=\", \"98\u00d755=\", \"85\u00d724=\"),\n  @(\"30\u00d732=\", \"47\u00d768=\", \"63\u00d721=\", \"91\u00d734=\", \"65\u00d791=\"),\n  @(\"95\u00d780=\", \"12\u00d775=\", \"89\u00d797=\", \"52\u00d714=\", \"36\u00d731=\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count -and $r -le $newValues.Count; $r++) {\n  $rowValues = $newValues[$r - 1]\n  for ($c = 1; $c -le $t.Columns.Count -and $c -le $rowValues.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $rowValues[$c - 1]\n  }\n}\n\n"}
